$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-20, -19.98, -34),
    @(-16, -15.95, -27),
    @(-12, -11.99, -20),
    @(-8, -8.0299999999999994, -14),
    @(-4, -3.97, -7),
    @(0, 0, 0),
    @(4, 4, 6),
    @(8, 8.01, 13),
    @(12, 12.03, 20),
    @(16, 15.95, 27),
    @(20, 19.96, 34)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 13).Value = $data[$i][0]
    $ws.Cells.Item($row, 14).Value = $data[$i][1]
    $ws.Cells.Item($row, 15).Value = $data[$i][2]
}

# Apply underline style (same as G6) to O6
$ws.Range("O6").Font.Underline = $true

# Q4 gets the same style applied but stays empty (no value)
$ws.Range("Q4").Font.Underline = $true

# Update selection to O6
$ws.Range("O6").Select()
